$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade row (row 4) to fix Bollinger Bands calculation
$ws.Range("A3:I3").Copy()
$ws.Range("A4:I4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 42635.800462962965
$ws.Range("B4").Value = $true
$ws.Range("C4").Value = 10085.129999999999
$ws.Range("D4").Value = 10020
$ws.Range("E4").Value = 82.03
$ws.Range("F4").Value = 80.9599
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = -1.3
$ws.Range("I4").Value = $false

# Column C's best-fit width grows slightly to fit the new, longer value (10085.13)
$ws.Columns.Item(3).ColumnWidth = 8.166666666666666
